$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-5 from 45170 to 45174
$ws.Range("C2").Value = 45174
$ws.Range("C3").Value = 45174
$ws.Range("C4").Value = 45174
$ws.Range("C5").Value = 45174
